# Apply the edits described by the commit:
# "toppings and decorations can be NULL for creation. created urls and views for cookiecreation"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Workbook window view (best effort - engine may not persist this, but it
#    is the semantically correct COM call for this change).
# ---------------------------------------------------------------------------
$win = $excel.ActiveWindow
$win.Left = 9840
$win.Top = -23960
$win.Width = 26780
$win.Height = 19680

# ---------------------------------------------------------------------------
# 2. New cell J4 - "googlesign_in." note next to the signup row.
#    (This is the first newly referenced shared string, so it must be set
#    first to land at shared-string index 38.)
# ---------------------------------------------------------------------------
$ws.Range("J4").Value = "googlesign_in.  "

# ---------------------------------------------------------------------------
# 3. New rows describing the "AllFlavors/AllCookieCutters/AllToppings/
#    AllDecorations" views, added under rows 17/19/21/23 (each of which
#    already held a single URL cell in column A).
#
#    Values are written column-by-column (rather than row-by-row) so that
#    newly introduced shared strings are interned in the same order as the
#    authored workbook (B17 "", C17.. view names, D17 the "get" call blurb,
#    then the E column "json: ..." responses).
# ---------------------------------------------------------------------------

# Column B - empty-json inputs marker ("")
$ws.Range("B17").Value = """"""
$ws.Range("B19").Value = """"""
$ws.Range("B21").Value = """"""
$ws.Range("B23").Value = """"""

# Column C - view names
$ws.Range("C17").Value = "AllFlavorsView"
$ws.Range("C19").Value = "AllCookieCuttersView"
$ws.Range("C21").Value = "AllToppingsView"
$ws.Range("C23").Value = "AllDecorationsView"

# Column D - model url / call description (wrapped, reuses style index 1)
$ws.Range("D17").Value = "get`n--auth header"
$ws.Range("D19").Value = "get`n--auth header"
$ws.Range("D21").Value = "get`n--auth header"
$ws.Range("D23").Value = "get`n--auth header"
$ws.Range("D17").WrapText = $true
$ws.Range("D19").WrapText = $true
$ws.Range("D21").WrapText = $true
$ws.Range("D23").WrapText = $true

# Column E - json responses
$ws.Range("E17").Value = "json: all flavors"
$ws.Range("E19").Value = "json: all cookie cutters"
$ws.Range("E21").Value = "json: all toppings"
$ws.Range("E23").Value = "json: all decorations"

# Column F - http status (reuses existing "200: ok" shared string)
$ws.Range("F17").Value = "200: ok"
$ws.Range("F19").Value = "200: ok"
$ws.Range("F21").Value = "200: ok"
$ws.Range("F23").Value = "200: ok"

# Row heights for the newly wrapped two-line cells.
$ws.Rows.Item(17).RowHeight = 32
$ws.Rows.Item(19).RowHeight = 32
$ws.Rows.Item(21).RowHeight = 32
$ws.Rows.Item(23).RowHeight = 32

# ---------------------------------------------------------------------------
# 4. Column C needs a best-fit width (to fit the new view-name strings).
#    Target stored width is 18.6640625; the engine quantizes ColumnWidth to
#    the nearest pixel bucket, so 17.83 is the closest achievable value.
# ---------------------------------------------------------------------------
$ws.Columns.Item(3).ColumnWidth = 17.83

# ---------------------------------------------------------------------------
# 5. Update the active selection to match the new edit location.
# ---------------------------------------------------------------------------
$ws.Range("H19").Select() | Out-Null
